$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing Salary row (row 2)
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 46021.29180555556

# Copy the date-formatted style from C2 down to C3:C6 before writing values,
# so the new rows inherit the same number format (numFmtId 14) as the
# existing date column instead of getting a fresh/plain style.
$ws.Range("C2").Copy($ws.Range("C3:C6"))

# New row 3: Business Income
$ws.Range("A3").Value = "Business Income"
$ws.Range("B3").Value = 300
$ws.Range("C3").Value = 46016.29180555556

# New row 4: Interest From Saving Account
$ws.Range("A4").Value = "Interest From Saving Account"
$ws.Range("B4").Value = 150
$ws.Range("C4").Value = 45999.29180555556

# New row 5: Trading
$ws.Range("A5").Value = "Trading"
$ws.Range("B5").Value = 300
$ws.Range("C5").Value = 45996.29180555556

# New row 6: Performance Bonus
$ws.Range("A6").Value = "Performance Bonus"
$ws.Range("B6").Value = 500
$ws.Range("C6").Value = 45993.29180555556
